# Trade #7 closed at 2026-02-17 13:08:12 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.58
$summary.Range("B4").Value = -0.42
$summary.Range("B5").Value = -1.2
$summary.Range("B6").Value = 7
$summary.Range("B8").Value = 5
$summary.Range("B9").Value = 28.57

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.58
$status.Range("D4").Value = 7
$status.Range("E4").Value = -0.42
$status.Range("F4").Value = -0.42
$status.Range("G4").Value = 28.57

# Write a string value into a cell without Excel's auto date/number
# detection silently reinterpreting it (and without leaving a stray
# NumberFormat override behind on the cell).
function Set-TextCell($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- Helper to append the new trade row (#7) to a trades sheet ---
function Add-TradeRow($ws) {
    $ws.Cells.Item(8, 1).Value = 7

    Set-TextCell $ws.Cells.Item(8, 2) "2026-02-17"
    Set-TextCell $ws.Cells.Item(8, 3) "13:08:06"
    Set-TextCell $ws.Cells.Item(8, 4) "MarketMaking"
    Set-TextCell $ws.Cells.Item(8, 5) "DOWN"

    $ws.Cells.Item(8, 6).Value = 0.785878
    $ws.Cells.Item(8, 7).Value = 0.78

    Set-TextCell $ws.Cells.Item(8, 8) "CLOSED"

    $ws.Cells.Item(8, 9).Value = -0.748
    $ws.Cells.Item(8, 10).Value = -0.01
    $ws.Cells.Item(8, 11).Value = 99.58
    $ws.Cells.Item(8, 12).Value = 0
    $ws.Cells.Item(8, 13).Value = 0
    $ws.Cells.Item(8, 14).Value = 0.6

    Set-TextCell $ws.Cells.Item(8, 15) "Normal spread capture: 19600 bps"
    Set-TextCell $ws.Cells.Item(8, 16) "early_exit"

    $ws.Cells.Item(8, 17).Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
